$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.325.79"
$ws.Cells.Item(2, 5).Value = "  -0.70%  "

$ws.Cells.Item(3, 4).Value = "1.785.92"
$ws.Cells.Item(3, 5).Value = "  -3.02%  "

$ws.Cells.Item(4, 5).Value = "  +0.01%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "224.55"
$ws.Cells.Item(5, 5).Value = "  -3.35%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.553"
$ws.Cells.Item(6, 5).Value = "  -3.07%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.999"
$ws.Cells.Item(7, 5).Value = "  -0.04%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "33.34"
$ws.Cells.Item(8, 5).Value = "  +4.72%  "

$ws.Cells.Item(9, 5).Value = "  -3.18%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0661"
$ws.Cells.Item(10, 5).Value = "  -3.88%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0931"
$ws.Cells.Item(11, 5).Value = "  -0.36%  "

$ws.Cells.Item(12, 4).Value = "2.042.15"
$ws.Cells.Item(12, 5).Value = "  -3.15%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "11.14"
$ws.Cells.Item(13, 5).Value = "  +7.40%  "

$ws.Cells.Item(14, 4).Value = "1.784.86"
$ws.Cells.Item(14, 5).Value = "  -3.16%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.633"
$ws.Cells.Item(15, 5).Value = "  -4.26%  "

$ws.Cells.Item(16, 4).Value = "34.280.73"
$ws.Cells.Item(16, 5).Value = "  -0.77%  "

$ws.Cells.Item(17, 5).Value = "  -2.44%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "69.02"
$ws.Cells.Item(18, 5).Value = "  -2.54%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "255.60"
$ws.Cells.Item(19, 5).Value = "  -1.76%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0744"
$ws.Cells.Item(20, 5).Value = "  -2.96%  "

$ws.Cells.Item(21, 5).Value = "  +0.10%  "

$ws.Cells.Item(22, 5).Value = "  -2.89%  "

$ws.Cells.Item(23, 5).Value = "  -5.23%  "

$ws.Cells.Item(24, 5).Value = "  -5.20%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "157.76"
$ws.Cells.Item(25, 5).Value = "  -0.92%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "16.42"
$ws.Cells.Item(26, 5).Value = "  -3.04%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.03"
$ws.Cells.Item(27, 5).Value = "  -2.89%  "

$ws.Cells.Item(28, 5).Value = "  -3.86%  "

$ws.Cells.Item(29, 5).Value = "  +0.01%  "

$ws.Cells.Item(30, 5).Value = "  -4.67%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.0515"
$ws.Cells.Item(31, 5).Value = "  -3.25%  "

$ws.Cells.Item(32, 5).Value = "  -2.86%  "

$ws.Cells.Item(33, 5).Value = "  -1.52%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.90"
$ws.Cells.Item(34, 5).Value = "  +3.98%  "

$ws.Cells.Item(35, 4).Value = "1.444.77"
$ws.Cells.Item(35, 5).Value = "  -7.63%  "

$ws.Cells.Item(36, 5).Value = "  -3.00%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0189"
$ws.Cells.Item(37, 5).Value = "  -2.15%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.625"
$ws.Cells.Item(38, 5).Value = "  -3.85%  "

$ws.Cells.Item(39, 2).Value = "Aave"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "83.08"
$ws.Cells.Item(39, 5).Value = "  -3.39%  "

$ws.Cells.Item(40, 2).Value = "MXToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.85"
$ws.Cells.Item(40, 5).Value = "  +0.26%  "

$ws.Cells.Item(41, 5).Value = "  -0.27%  "

$ws.Cells.Item(42, 5).Value = "  -3.95%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.07"
$ws.Cells.Item(43, 5).Value = "  -4.00%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0506"
$ws.Cells.Item(44, 5).Value = "  -4.35%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.06"
$ws.Cells.Item(45, 5).Value = "  -2.21%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "12.43"
$ws.Cells.Item(46, 5).Value = "  -0.04%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "5.83"
$ws.Cells.Item(47, 5).Value = "  -0.60%  "

$ws.Cells.Item(48, 4).Value = "1.941.97"
$ws.Cells.Item(48, 5).Value = "  -2.87%  "

$ws.Cells.Item(49, 5).Value = "  -0.03%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "98.30"
$ws.Cells.Item(50, 5).Value = "  -2.00%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "49.90"
$ws.Cells.Item(51, 5).Value = "  -3.89%  "

Write-Output "applied changes"